# "Implemented delete function for serialisation"
#
# The FOMS menu_list.xlsx sheet holds one menu item per row (id, name,
# price, branch, category, description) starting at row 2 (row 1 is the
# header). Exercising the app's new "delete menu item" function removed
# the "pepsi" row (row 9: 61606264-3573-4492-abed-5a026cdfc717 / pepsi /
# 2.1 / JE / drink / Drink) from the data.
#
# Deleting the entire row (rather than just clearing its contents) shifts
# every row below it up by one, which is why the sheet's used range goes
# from A1:F76 down to A1:F75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(9).Delete()
